$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 461.50722525496786
$ws.Range("C2").Value = 813.75693527432168
$ws.Range("D2").Value = 417.24887536970186
$ws.Range("E2").Value = 736.47082127557076

# Row 3 values
$ws.Range("B3").Value = 444.87390981478114
$ws.Range("C3").Value = 689.23854332370718
$ws.Range("D3").Value = 641.6741646289305
$ws.Range("E3").Value = 476.2203269271044

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
